$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet's displayed tab title (sheet name in workbook.xml)
$ws.Name = "UniformA"

# Append the new row of data (row 16), following the same pattern as the
# preceding rows (A = index, B = shared text label, C:P = 1)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}

# Match the bold/bordered/centered style used by the other cells in column A
# by copying the formatting from the cell directly above it.
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
